$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy date-header formatting from AY1:BA1 to BB1:BD1 ---
$ws.Range("AY1:BA1").Copy()
$ws.Range("BB1:BD1").PasteSpecial(-4122)
$ws.Range("BB1").Value = "31/12/2023"
$ws.Range("BC1").Value = "31/03/2024"
$ws.Range("BD1").Value = "30/06/2024"

# --- Data rows: BB/BC/BD columns, rows 2-80 ---
# Row 2
$ws.Range("BB2").Value = 503259.008
$ws.Range("BC2").Value = 430204
$ws.Range("BD2").Value = 449380.992
# Row 3
$ws.Range("BB3").Value = 290348.992
$ws.Range("BC3").Value = 211100.992
$ws.Range("BD3").Value = 232900
# Row 4
$ws.Range("BB4").Value = 146200
$ws.Range("BC4").Value = 92979
$ws.Range("BD4").Value = 85921
# Row 5
$ws.Range("BB5").Value = 2327
$ws.Range("BC5").Value = 2385
$ws.Range("BD5").Value = 2437
# Row 6
$ws.Range("BB6").Value = 109470
$ws.Range("BC6").Value = 91445
$ws.Range("BD6").Value = 128225
# Row 7
$ws.Range("BB7").Value = 1945
$ws.Range("BC7").Value = 1862
$ws.Range("BD7").Value = 1634
# Row 8
$ws.Range("BB8").Value = 0
$ws.Range("BC8").Value = 0
$ws.Range("BD8").Value = 0
# Row 9
$ws.Range("BB9").Value = 0
$ws.Range("BC9").Value = 0
$ws.Range("BD9").Value = 0
# Row 10
$ws.Range("BB10").Value = 30407
$ws.Range("BC10").Value = 22316
$ws.Range("BD10").Value = 9247
# Row 11
$ws.Range("BB11").Value = 0
$ws.Range("BC11").Value = 114
$ws.Range("BD11").Value = 5436
# Row 12
$ws.Range("BB12").Value = 105547
$ws.Range("BC12").Value = 110828
$ws.Range("BD12").Value = 108102
# Row 13
$ws.Range("BB13").Value = 776
$ws.Range("BC13").Value = 795
$ws.Range("BD13").Value = 812
# Row 14
$ws.Range("BB14").Value = 0
$ws.Range("BC14").Value = 0
$ws.Range("BD14").Value = 0
# Row 15
$ws.Range("BB15").Value = 0
$ws.Range("BC15").Value = 0
$ws.Range("BD15").Value = 0
# Row 16
$ws.Range("BB16").Value = 32234
$ws.Range("BC16").Value = 32234
$ws.Range("BD16").Value = 32234
# Row 17
$ws.Range("BB17").Value = 0
$ws.Range("BC17").Value = 0
$ws.Range("BD17").Value = 0
# Row 18
$ws.Range("BB18").Value = 0
$ws.Range("BC18").Value = 0
$ws.Range("BD18").Value = 0
# Row 19
$ws.Range("BB19").Value = 50069
$ws.Range("BC19").Value = 49598
$ws.Range("BD19").Value = 49940
# Row 20
$ws.Range("BB20").Value = 487
$ws.Range("BC20").Value = 487
$ws.Range("BD20").Value = 487
# Row 21
$ws.Range("BB21").Value = 11216
$ws.Range("BC21").Value = 16682
$ws.Range("BD21").Value = 13649
# Row 22
$ws.Range("BB22").Value = 0
$ws.Range("BC22").Value = 0
$ws.Range("BD22").Value = 0
# Row 23
$ws.Range("BB23").Value = 19517
$ws.Range("BC23").Value = 20180
$ws.Range("BD23").Value = 20195
# Row 24
$ws.Range("BB24").Value = 87846
$ws.Range("BC24").Value = 88095
$ws.Range("BD24").Value = 88184
# Row 25
$ws.Range("BB25").Value = 0
$ws.Range("BC25").Value = 0
$ws.Range("BD25").Value = 0
# Row 26
$ws.Range("BB26").Value = 503259.008
$ws.Range("BC26").Value = 430204
$ws.Range("BD26").Value = 449380.992
# Row 27
$ws.Range("BB27").Value = 238220.992
$ws.Range("BC27").Value = 173974
$ws.Range("BD27").Value = 205596.992
# Row 28
$ws.Range("BB28").Value = 4792
$ws.Range("BC28").Value = 4744
$ws.Range("BD28").Value = 2113
# Row 29
$ws.Range("BB29").Value = 145926
$ws.Range("BC29").Value = 83575
$ws.Range("BD29").Value = 99899
# Row 30
$ws.Range("BB30").Value = 13418
$ws.Range("BC30").Value = 11593
$ws.Range("BD30").Value = 12189
# Row 31
$ws.Range("BB31").Value = 42172
$ws.Range("BC31").Value = 44167
$ws.Range("BD31").Value = 68724
# Row 32
$ws.Range("BB32").Value = 0
$ws.Range("BC32").Value = 0
$ws.Range("BD32").Value = 0
# Row 33
$ws.Range("BB33").Value = 0
$ws.Range("BC33").Value = 0
$ws.Range("BD33").Value = 0
# Row 34
$ws.Range("BB34").Value = 31913
$ws.Range("BC34").Value = 29895
$ws.Range("BD34").Value = 22672
# Row 35
$ws.Range("BB35").Value = 0
$ws.Range("BC35").Value = 0
$ws.Range("BD35").Value = 0
# Row 36
$ws.Range("BB36").Value = 0
$ws.Range("BC36").Value = 0
$ws.Range("BD36").Value = 0
# Row 37
$ws.Range("BB37").Value = 84725
$ws.Range("BC37").Value = 75479
$ws.Range("BD37").Value = 75145
# Row 38
$ws.Range("BB38").Value = 13140
$ws.Range("BC38").Value = 12475
$ws.Range("BD38").Value = 11824
# Row 39
$ws.Range("BB39").Value = 0
$ws.Range("BC39").Value = 0
$ws.Range("BD39").Value = 0
# Row 40
$ws.Range("BB40").Value = 22033
$ws.Range("BC40").Value = 21766
$ws.Range("BD40").Value = 21309
# Row 41
$ws.Range("BB41").Value = 2258
$ws.Range("BC41").Value = 2794
$ws.Range("BD41").Value = 1948
# Row 42
$ws.Range("BB42").Value = 0
$ws.Range("BC42").Value = 0
$ws.Range("BD42").Value = 0
# Row 43
$ws.Range("BB43").Value = 47294
$ws.Range("BC43").Value = 38444
$ws.Range("BD43").Value = 40064
# Row 44
$ws.Range("BB44").Value = 0
$ws.Range("BC44").Value = 0
$ws.Range("BD44").Value = 0
# Row 45
$ws.Range("BB45").Value = 0
$ws.Range("BC45").Value = 0
$ws.Range("BD45").Value = 0
# Row 46
$ws.Range("BB46").Value = 3278
$ws.Range("BC46").Value = 3316
$ws.Range("BD46").Value = 4056
# Row 47
$ws.Range("BB47").Value = 177034.992
$ws.Range("BC47").Value = 177435.008
$ws.Range("BD47").Value = 164583.008
# Row 48
$ws.Range("BB48").Value = 243022
$ws.Range("BC48").Value = 243022
$ws.Range("BD48").Value = 243022
# Row 49
$ws.Range("BB49").Value = 0
$ws.Range("BC49").Value = 0
$ws.Range("BD49").Value = 0
# Row 50
$ws.Range("BB50").Value = 154
$ws.Range("BC50").Value = 150
$ws.Range("BD50").Value = 145
# Row 51
$ws.Range("BB51").Value = -9665
$ws.Range("BC51").Value = -9665
$ws.Range("BD51").Value = -9665
# Row 52
$ws.Range("BB52").Value = -115044
$ws.Range("BC52").Value = -118395
$ws.Range("BD52").Value = -144480
# Row 53
$ws.Range("BB53").Value = 0
$ws.Range("BC53").Value = 0
$ws.Range("BD53").Value = 0
# Row 54
$ws.Range("BB54").Value = 0
$ws.Range("BC54").Value = 0
$ws.Range("BD54").Value = 0
# Row 55
$ws.Range("BB55").Value = 58568
$ws.Range("BC55").Value = 62323
$ws.Range("BD55").Value = 75561
# Row 56
$ws.Range("BB56").Value = 0
$ws.Range("BC56").Value = 0
$ws.Range("BD56").Value = 0
# Row 59
$ws.Range("BB59").Value = 339667.008
$ws.Range("BC59").Value = 34645
$ws.Range("BD59").Value = 55129
# Row 60
$ws.Range("BB60").Value = -294778.976
$ws.Range("BC60").Value = -35448
$ws.Range("BD60").Value = -57429
# Row 61
$ws.Range("BB61").Value = 44888
$ws.Range("BC61").Value = -803
$ws.Range("BD61").Value = -2300
# Row 62
$ws.Range("BB62").Value = 9
$ws.Range("BC62").Value = -241
$ws.Range("BD62").Value = -563
# Row 63
$ws.Range("BB63").Value = -15266
$ws.Range("BC63").Value = -9734
$ws.Range("BD63").Value = -8905
# Row 64
$ws.Range("BB64").Value = 0
$ws.Range("BC64").Value = 0
$ws.Range("BD64").Value = 0
# Row 65
$ws.Range("BB65").Value = -149
$ws.Range("BC65").Value = 9575
$ws.Range("BD65").Value = 3066
# Row 66
$ws.Range("BB66").Value = -1400
$ws.Range("BC66").Value = 0
$ws.Range("BD66").Value = 0
# Row 67
$ws.Range("BB67").Value = -389
$ws.Range("BC67").Value = 0
$ws.Range("BD67").Value = 0
# Row 68
$ws.Range("BB68").Value = -1686
$ws.Range("BC68").Value = 513
$ws.Range("BD68").Value = -15966
# Row 69
$ws.Range("BB69").Value = 616
$ws.Range("BC69").Value = 3549
$ws.Range("BD69").Value = 4489
# Row 70
$ws.Range("BB70").Value = -2302
$ws.Range("BC70").Value = -3036
$ws.Range("BD70").Value = -20455
# Row 74
$ws.Range("BB74").Value = 26007
$ws.Range("BC74").Value = -690
$ws.Range("BD74").Value = -24668
# Row 75
$ws.Range("BB75").Value = -2565
$ws.Range("BC75").Value = -2709
$ws.Range("BD75").Value = -1340
# Row 76
$ws.Range("BB76").Value = 813
$ws.Range("BC76").Value = 133
$ws.Range("BD76").Value = 1099
# Row 79
$ws.Range("BB79").Value = 1722
$ws.Range("BC79").Value = -89
$ws.Range("BD79").Value = -1181
# Row 80
$ws.Range("BB80").Value = 25977
$ws.Range("BC80").Value = -3355
$ws.Range("BD80").Value = -26090

# --- Blank placeholder rows (no numeric data; extend row formatting out to BD) ---
# These rows are section separators whose A-E etc. columns are also blank placeholders;
# replicate that by copying the blank format from column BA into BB:BD.
$ws.Range("BA57").Copy()
$ws.Range("BB57:BD57").PasteSpecial(-4122)
$ws.Range("BA58").Copy()
$ws.Range("BB58:BD58").PasteSpecial(-4122)
$ws.Range("BA71").Copy()
$ws.Range("BB71:BD71").PasteSpecial(-4122)
$ws.Range("BA72").Copy()
$ws.Range("BB72:BD72").PasteSpecial(-4122)
$ws.Range("BA73").Copy()
$ws.Range("BB73:BD73").PasteSpecial(-4122)
$ws.Range("BA77").Copy()
$ws.Range("BB77:BD77").PasteSpecial(-4122)
$ws.Range("BA78").Copy()
$ws.Range("BB78:BD78").PasteSpecial(-4122)
